# Normalise separators used in the "research direction" / "status" free-text
# cells: semicolons (and one stray full-width comma) become plain ", ",
# and a couple of stray/trailing spaces & a missing capital are tidied up.
# Only Sheet1 is involved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = "Netease Fuxi AI Lab, Alibaba"
$ws.Range("I5").Value = "Postdoctoral Fellow of UdeM/MILA, Associate Reseacher, Tianjin University"
$ws.Range("G6").Value = "Reinforcement Learning, Transfer Learning, Multiagent Learning"
$ws.Range("I6").Value = "Postdoc at University of Alberta, Associate Professor, Nanjing University"
$ws.Range("G7").Value = "Multiagent Systems, Deep Reinforcement Learning, Evolutionary Algorithm"
$ws.Range("G8").Value = "Reinforcement Learning, Multiagent Reinforcement Learning"
$ws.Range("G9").Value = "Model based RL, Diffusion for RL, LLM"
$ws.Range("I9").Value = "Postdoctoral Fellow, Imperial College London"

# Restore the view: the active window had scrolled/selected down to I12;
# bring it back so I10 is selected and column C is the left-most visible
# column again.
$ws.Range("I10").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
